$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Simple Format")
$ws.Activate()

# A4: "Big font" at size 72
$ws.Range("A4").Value = "Big font"
$ws.Range("A4").Font.Size = 72

# A5: "Medium font" at size 36
$ws.Range("A5").Value = "Medium font"
$ws.Range("A5").Font.Size = 36

# A6: rich text with varying font sizes within the same cell
$ws.Range("A6").Value = "Big, medium, and small fonts."
$ws.Range("A6").Characters(1, 3).Font.Size = 72
$ws.Range("A6").Characters(6, 6).Font.Size = 36
$ws.Range("A6").Characters(18, 5).Font.Size = 8

$ws.Rows.Item(4).RowHeight = 92.25
$ws.Rows.Item(5).RowHeight = 46.5
$ws.Rows.Item(6).RowHeight = 92.25

$ws.Range("A6").Select()
